$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Usersite Menu - Documents" row right after row 7
# (row 9 stays blank as a separator before the row-10 section header)
$ws.Range("B8").Value = "Usersite Menu - Documents"
$ws.Range("C8").Value = "xpath"
$ws.Range("D8").Value = "id"
$ws.Range("E8").Value = "Documents"

# Update the visible selection to the newly added row (also clears the
# previous topLeftCell scroll position since A1 is back in view)
$ws.Range("B8").Select() | Out-Null
